# Update PPT per Will's email
# - Slide 3: retitle "KennySync Must:" -> "Mandatory Goals"
# - Slide 5 (Team Organization): replace Eric/Will/Tim placeholders with new bullet content
# - Slide 6 (Timeline): flesh out milestone bullets with due-dates and sub-bullets
# - New slide inserted after Timeline: "Background literature"
# - Old "Questions?" slide (now after the new slide): drop "Papers, papers, papers" bullet

$p = $ppt.ActivePresentation

# --- Slide 3: "KennySync Must:" -> "Mandatory Goals" ---
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Mandatory Goals"

# --- Slide 5: Team Organization content ---
$s5 = $p.Slides.Item(5)
$tf5 = $s5.Shapes.Item(2).TextFrame
$tf5.TextRange.Text = "Minimal`rHighly ad-hoc: team works in one- or two-week sprints, completing significant features each week`rUsing Ruby 1.9; Tim may serve as domain expert in language`rFeatures, requirements, & planning decided by consensus`rTeam of 3 is small enough to make this work`r"
$tf5.TextRange.Paragraphs(5,1).IndentLevel = 2

# --- Slide 6: Timeline content ---
$s6 = $p.Slides.Item(6)
$tf6 = $s6.Shapes.Item(2).TextFrame
$tf6.TextRange.Text = "Week 5 Day 4 – Milestone 2 due (Literature review and design)`rTeam will have read and understand algorithm, variations`rWeek 6 Day 4 – Have connected Ruby nodes`rWeek 8 Day 4 – Have majority of algorithm implemented`rWeek 10 Day 4 – Milestone 3 due (Project delivery, Demo, Report, Presentation)`rTeam will have implemented visualization for demo to class"
$tf6.TextRange.Paragraphs(2,1).IndentLevel = 2
$tf6.TextRange.Paragraphs(6,1).IndentLevel = 2

# --- New slide: duplicate the trailing "Questions?" slide, then repurpose the ---
# --- original as "Background literature" (keeps the duplicate as the new    ---
# --- final "Questions?" slide, matching the inserted-before-last ordering). ---
$s7 = $p.Slides.Item(7)
$dup = $s7.Duplicate()
$newQuestions = $dup.Item(1)

$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Background literature"
$tf7 = $s7.Shapes.Item(2).TextFrame
$tf7.TextRange.Text = "Leslie Lamport wrote it all`rFirst description of Paxos, couched in Greek governance terms`r“Simple” Paxos description followed as brief (13-page) note`rThree other variations:`rCheap`rFast`rGeneralized"
$tf7.TextRange.Paragraphs(2,1).IndentLevel = 2
$tf7.TextRange.Paragraphs(3,1).IndentLevel = 2
$tf7.TextRange.Paragraphs(4,1).IndentLevel = 2
$tf7.TextRange.Paragraphs(5,1).IndentLevel = 3
$tf7.TextRange.Paragraphs(6,1).IndentLevel = 3
$tf7.TextRange.Paragraphs(7,1).IndentLevel = 3

# --- Final "Questions?" slide (the duplicate): drop the "Papers, papers, papers" bullet ---
$newQuestions.Shapes.Item(2).TextFrame.TextRange.Text = ""

Write-Output "Slide count: $($p.Slides.Count)"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    Write-Output "$i : $($p.Slides.Item($i).Shapes.Item(1).TextFrame.TextRange.Text)"
}
